$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
  "G2" = 32.74870266666667
  "H2" = 98.24610799999999
  "I2" = 0.6379427096830631
  "J2" = 0.6379427096830631
  "M2" = 43.20983933333334
  "N2" = 129.629518
  "O2" = 0.6315174248499266
  "P2" = 0.6315174248499266
  "Q2" = 1415.066180601772
  "R2" = 12735.59562541594
  "S2" = 0.4028719372208323
  "T2" = 0.4028719372208323
  "G3" = 32.74870266666667
  "H3" = 98.24610799999999
  "I3" = 0.6379427096830631
  "J3" = 0.6379427096830631
  "O3" = 0.01329126241460306
  "P3" = 0.01329126241460306
  "Q3" = 29.78225968171466
  "R3" = 268.040337135432
  "S3" = 0.008479063959880529
  "T3" = 0.008479063959880529
  "G4" = 32.74870266666667
  "H4" = 98.24610799999999
  "I4" = 0.6379427096830631
  "J4" = 0.6379427096830631
  "M4" = 24.10288866666667
  "N4" = 72.308666
  "O4" = 0.3522668544262691
  "P4" = 0.3522668544262691
  "Q4" = 789.3383343524365
  "R4" = 7104.045009171928
  "S4" = 0.2247260716442233
  "T4" = 0.2247260716442232
  "G5" = 32.74870266666667
  "H5" = 98.24610799999999
  "I5" = 0.6379427096830631
  "J5" = 0.6379427096830631
  "M5" = 0.200098
  "N5" = 0.600294
  "O5" = 0.002924458309201317
  "P5" = 0.002924458309201317
  "Q5" = 6.552949906194667
  "R5" = 58.976549155752
  "S5" = 0.001865636858127037
  "T5" = 0.001865636858127037
  "I6" = 0.01497516366297488
  "J6" = 0.01497516366297488
  "M6" = 43.20983933333334
  "N6" = 129.629518
  "O6" = 0.6315174248499266
  "P6" = 0.6315174248499266
  "Q6" = 33.21747756782134
  "R6" = 298.9572981103921
  "S6" = 0.009457076793148092
  "T6" = 0.009457076793148092
  "I7" = 0.01497516366297488
  "J7" = 0.01497516366297488
  "O7" = 0.01329126241460306
  "P7" = 0.01329126241460306
  "S7" = 0.0001990388299462276
  "T7" = 0.0001990388299462276
  "I8" = 0.01497516366297488
  "J8" = 0.01497516366297488
  "M8" = 24.10288866666667
  "N8" = 72.308666
  "O8" = 0.3522668544262691
  "P8" = 0.3522668544262691
  "Q8" = 18.52904745672267
  "R8" = 166.761427110504
  "S8" = 0.005275253798074729
  "T8" = 0.005275253798074728
  "I9" = 0.01497516366297488
  "J9" = 0.01497516366297488
  "M9" = 0.200098
  "N9" = 0.600294
  "O9" = 0.002924458309201317
  "P9" = 0.002924458309201317
  "Q9" = 0.153824937304
  "R9" = 1.384424435736
  "S9" = 0.00004379424180583654
  "T9" = 0.00004379424180583654
  "G10" = 17.70628
  "H10" = 53.11884
  "I10" = 0.3449172431830183
  "J10" = 0.3449172431830183
  "M10" = 43.20983933333334
  "N10" = 129.629518
  "O10" = 0.6315174248499266
  "P10" = 0.6315174248499266
  "Q10" = 765.0855139910135
  "R10" = 6885.769625919121
  "S10" = 0.2178212492012756
  "T10" = 0.2178212492012756
  "G11" = 17.70628
  "H11" = 53.11884
  "I11" = 0.3449172431830183
  "J11" = 0.3449172431830183
  "O11" = 0.01329126241460306
  "P11" = 0.01329126241460306
  "Q11" = 16.10240974504
  "R11" = 144.92168770536
  "S11" = 0.004584385590466955
  "T11" = 0.004584385590466956
  "G12" = 17.70628
  "H12" = 53.11884
  "I12" = 0.3449172431830183
  "J12" = 0.3449172431830183
  "M12" = 24.10288866666667
  "N12" = 72.308666
  "O12" = 0.3522668544262691
  "P12" = 0.3522668544262691
  "Q12" = 426.7724955408267
  "R12" = 3840.95245986744
  "S12" = 0.1215029122934624
  "T12" = 0.1215029122934624
  "G13" = 17.70628
  "H13" = 53.11884
  "I13" = 0.3449172431830183
  "J13" = 0.3449172431830183
  "M13" = 0.200098
  "N13" = 0.600294
  "O13" = 0.002924458309201317
  "P13" = 0.002924458309201317
  "Q13" = 3.54299121544
  "R13" = 31.88692093896
  "S13" = 0.001008696097813389
  "T13" = 0.001008696097813389
  "G14" = 0.111134
  "H14" = 0.333402
  "I14" = 0.00216488347094373
  "J14" = 0.00216488347094373
  "M14" = 43.20983933333334
  "N14" = 129.629518
  "O14" = 0.6315174248499266
  "P14" = 0.6315174248499266
  "Q14" = 4.802082284470668
  "R14" = 43.218740560236
  "S14" = 0.001367161634670555
  "T14" = 0.001367161634670555
  "G15" = 0.111134
  "H15" = 0.333402
  "I15" = 0.00216488347094373
  "J15" = 0.00216488347094373
  "O15" = 0.01329126241460306
  "P15" = 0.01329126241460306
  "Q15" = 0.101067260012
  "R15" = 0.9096053401079999
  "S15" = 0.00002877403430934982
  "T15" = 0.00002877403430934982
  "G16" = 0.111134
  "H16" = 0.333402
  "I16" = 0.00216488347094373
  "J16" = 0.00216488347094373
  "M16" = 24.10288866666667
  "N16" = 72.308666
  "O16" = 0.3522668544262691
  "P16" = 0.3522668544262691
  "Q16" = 2.678650429081333
  "R16" = 24.107853861732
  "S16" = 0.0007626166905087714
  "T16" = 0.0007626166905087713
  "G17" = 0.111134
  "H17" = 0.333402
  "I17" = 0.00216488347094373
  "J17" = 0.00216488347094373
  "M17" = 0.200098
  "N17" = 0.600294
  "O17" = 0.002924458309201317
  "P17" = 0.002924458309201317
  "Q17" = 0.022237691132
  "R17" = 0.200139220188
  "S17" = 0.000006331111455053981
  "T17" = 0.000006331111455053981
}

foreach ($addr in $updates.Keys) {
  $ws.Range($addr).Value = $updates[$addr]
}
